$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)
$shape = $s.Shapes.Item(4)
$table = $shape.Table
$cell = $table.Cell(3, 2)
$cell.Shape.TextFrame.TextRange.Text = "1810.52"
